$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7112
$ws1.Range("F4").Value = 465
$ws1.Range("F7").Value = 163
$ws1.Range("F13").Value = 451
$ws1.Range("F15").Value = 1837
$ws1.Range("F17").Value = 3683
$ws1.Range("F18").Value = 27
$ws1.Range("F20").Value = 84
$ws1.Range("F21").Value = 29
$ws1.Range("F23").Value = 2320
$ws1.Range("F25").Value = 272
$ws1.Range("F26").Value = 11
$ws1.Range("F32").Value = 1357
$ws1.Range("F33").Value = 122

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7112
$ws4.Range("F4").Value = 465
$ws4.Range("F8").Value = 163
$ws4.Range("F14").Value = 451
$ws4.Range("F16").Value = 1837
$ws4.Range("F18").Value = 3683
$ws4.Range("F19").Value = 27
$ws4.Range("F21").Value = 84
$ws4.Range("F22").Value = 29
$ws4.Range("F24").Value = 2320
$ws4.Range("F26").Value = 272
$ws4.Range("F27").Value = 11
$ws4.Range("F33").Value = 1357
$ws4.Range("F34").Value = 122

$wb.Save()
